$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 118-138 (columns A-T), per the commit diff:
# the Kiwi price rows for Terminal Hortofruticola Agro Chillan were
# re-sequenced (Fecha/Volumen/Precio/Unidad/Precio-$-Kg values updated)
# and two brand-new rows (137-138, Fecha serial 44400 = 2021-07-23) were
# appended at the end, growing the sheet from A1:T136 to A1:T138.
$data = @(
  @(118, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44449, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 100, 12000, 12500, 12250, '$/bandeja 18 kilos', 'Provincia de Curicó', 681, 18),
  @(119, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44449, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 60, 11000, 11500, 11250, '$/bandeja 18 kilos', 'Provincia de Curicó', 625, 18),
  @(120, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44438, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 12000, 12500, 12250, '$/bandeja 18 kilos', 'Provincia de Curicó', 681, 18),
  @(121, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44438, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 60, 11000, 11500, 11250, '$/bandeja 18 kilos', 'Provincia de Curicó', 625, 18),
  @(122, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44442, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 12000, 12500, 12250, '$/bandeja 18 kilos', 'Provincia de Curicó', 681, 18),
  @(123, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44442, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 100, 11000, 11500, 11250, '$/bandeja 18 kilos', 'Provincia de Curicó', 625, 18),
  @(124, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44435, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 560, 12000, 12500, 12250, '$/bandeja 18 kilos', 'Provincia de Curicó', 681, 18),
  @(125, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44435, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 400, 11000, 11500, 11250, '$/bandeja 18 kilos', 'Provincia de Curicó', 625, 18),
  @(126, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44319, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 160, 10000, 11000, 10500, '$/bandeja 10 kilos', 'Provincia de Curicó', 1050, 10),
  @(127, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44319, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 120, 8500, 9000, 8750, '$/bandeja 10 kilos', 'Provincia de Curicó', 875, 10),
  @(128, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44376, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 10000, 11000, 10500, '$/bandeja 18 kilos', 'Provincia de Curicó', 583, 18),
  @(129, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44376, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 120, 8500, 9000, 8750, '$/bandeja 18 kilos', 'Provincia de Curicó', 486, 18),
  @(130, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44412, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 8000, 8500, 8250, '$/bandeja 10 kilos', 'Provincia de Curicó', 825, 10),
  @(131, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44412, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 80, 7000, 7500, 7250, '$/bandeja 10 kilos', 'Provincia de Curicó', 725, 10),
  @(132, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44314, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 8500, 9000, 8750, '$/bandeja 10 kilos', 'Provincia de Curicó', 875, 10),
  @(133, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44448, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 12000, 12500, 12250, '$/bandeja 18 kilos', 'Provincia de Curicó', 681, 18),
  @(134, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44448, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 60, 11000, 11500, 11250, '$/bandeja 18 kilos', 'Provincia de Curicó', 625, 18),
  @(135, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44399, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 10000, 11000, 10500, '$/bandeja 18 kilos', 'Provincia de Curicó', 583, 18),
  @(136, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44399, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 80, 9000, 9000, 9000, '$/bandeja 18 kilos', 'Provincia de Curicó', 500, 18),
  @(137, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44400, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Primera', 120, 11000, 12000, 11500, '$/bandeja 18 kilos', 'Provincia de Curicó', 639, 18),
  @(138, 7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44400, 16, 'Fruta', 100101, 'Berries', 100101007, 'Kiwi', 'Hayward', 'Segunda', 120, 9000, 10000, 9500, '$/bandeja 18 kilos', 'Provincia de Curicó', 528, 18)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($col = 1; $col -le 20; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}

# Rows 137-138 are brand new; give the Fecha (date) column the same
# date number format used by the rest of the column (D2:D136).
$dateFmt = $ws.Cells.Item(136, 4).NumberFormat
$ws.Cells.Item(137, 4).NumberFormat = $dateFmt
$ws.Cells.Item(138, 4).NumberFormat = $dateFmt
